$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the existing "sum" header (G1) to the new "Save" header (H1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text and value
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
